$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove Belle Ho from the roster (row 5); everyone below shifts up one row.
$ws.Rows(5).Delete()

# Add "Team" / "Dataset" header cells, matching the bold header formatting
# already used by the other header cells in row 1.
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Team"
$ws.Range("E1").Value = "Dataset"

# Populate each member's assigned Team / Dataset.
$ws.Range("D2").Value = "Pitching"
$ws.Range("E2").Value = "Relief Pitching"

$ws.Range("D3").Value = "Hitting"
$ws.Range("E3").Value = "Standard Hitting"

$ws.Range("D4").Value = "Hitting"
$ws.Range("E4").Value = "Advanced Hitting"

$ws.Range("D5").Value = "Fielding"
$ws.Range("E5").Value = "Standard Fielding"

$ws.Range("D6").Value = "Pitching"
$ws.Range("E6").Value = "Standard Pitching"

$ws.Range("D7").Value = "Hitting"
$ws.Range("E7").Value = "Sabermetric Hitting"

$ws.Range("D8").Value = "Fielding"
$ws.Range("E8").Value = "Outfield Fielding"

$ws.Range("D9").Value = "Pitching"
$ws.Range("E9").Value = "Starting Pitching"

# Match the existing shaded formatting in column C (rows 5-6) by applying the
# same format to the Dataset cells in those rows.
$ws.Range("C5").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = "Standard Fielding"

$ws.Range("C6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = "Standard Pitching"
